$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.95578266666667
$ws.Range("H2").Value = 59.867348
$ws.Range("I2").Value = 0.0117373419656925
$ws.Range("J2").Value = 0.0117373419656925
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.135193666666667
$ws.Range("N2").Value = 3.405581
$ws.Range("O2").Value = 0.153770120695047
$ws.Range("P2").Value = 0.153770120695047
$ws.Range("Q2").Value = 22.65367809657645
$ws.Range("R2").Value = 203.883102869188
$ws.Range("S2").Value = 0.001804852490703575
$ws.Range("T2").Value = 0.001804852490703575
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.95578266666667
$ws.Range("H3").Value = 59.867348
$ws.Range("I3").Value = 0.0117373419656925
$ws.Range("J3").Value = 0.0117373419656925
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.287366
$ws.Range("N3").Value = 6.862098
$ws.Range("O3").Value = 0.3098401235152652
$ws.Range("P3").Value = 0.3098401235152652
$ws.Range("Q3").Value = 45.64617877512267
$ws.Range("R3").Value = 410.815608976104
$ws.Range("S3").Value = 0.003636699484391069
$ws.Range("T3").Value = 0.003636699484391069
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.95578266666667
$ws.Range("H4").Value = 59.867348
$ws.Range("I4").Value = 0.0117373419656925
$ws.Range("J4").Value = 0.0117373419656925
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.138589666666667
$ws.Range("N4").Value = 9.415769000000001
$ws.Range("O4").Value = 0.425144471843918
$ws.Range("P4").Value = 0.425144471843918
$ws.Range("Q4").Value = 62.63301326784578
$ws.Range("R4").Value = 563.697119410612
$ws.Range("S4").Value = 0.004990066050855791
$ws.Range("T4").Value = 0.004990066050855791
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.95578266666667
$ws.Range("H5").Value = 59.867348
$ws.Range("I5").Value = 0.0117373419656925
$ws.Range("J5").Value = 0.0117373419656925
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8212579999999999
$ws.Range("N5").Value = 2.463774
$ws.Range("O5").Value = 0.1112452839457698
$ws.Range("P5").Value = 0.1112452839457698
$ws.Range("Q5").Value = 16.38884616126133
$ws.Range("R5").Value = 147.499615451352
$ws.Range("S5").Value = 0.001305723939742062
$ws.Range("T5").Value = 0.001305723939742062
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1637.343343333333
$ws.Range("H6").Value = 4912.03003
$ws.Range("I6").Value = 0.9630320723052701
$ws.Range("J6").Value = 0.9630320723052702
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.135193666666667
$ws.Range("N6").Value = 3.405581
$ws.Range("O6").Value = 0.153770120695047
$ws.Range("P6").Value = 0.153770120695047
$ws.Range("Q6").Value = 1858.701793510825
$ws.Range("R6").Value = 16728.31614159743
$ws.Range("S6").Value = 0.1480855579915826
$ws.Range("T6").Value = 0.1480855579915826
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1637.343343333333
$ws.Range("H7").Value = 4912.03003
$ws.Range("I7").Value = 0.9630320723052701
$ws.Range("J7").Value = 0.9630320723052702
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.287366
$ws.Range("N7").Value = 6.862098
$ws.Range("O7").Value = 0.3098401235152652
$ws.Range("P7").Value = 0.3098401235152652
$ws.Range("Q7").Value = 3745.203493866993
$ws.Range("R7").Value = 33706.83144480293
$ws.Range("S7").Value = 0.2983859762322267
$ws.Range("T7").Value = 0.2983859762322268
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1637.343343333333
$ws.Range("H8").Value = 4912.03003
$ws.Range("I8").Value = 0.9630320723052701
$ws.Range("J8").Value = 0.9630320723052702
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.138589666666667
$ws.Range("N8").Value = 9.415769000000001
$ws.Range("O8").Value = 0.425144471843918
$ws.Range("P8").Value = 0.425144471843918
$ws.Range("Q8").Value = 5138.948898171452
$ws.Range("R8").Value = 46250.54008354307
$ws.Range("S8").Value = 0.4094277617489779
$ws.Range("T8").Value = 0.409427761748978
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1637.343343333333
$ws.Range("H9").Value = 4912.03003
$ws.Range("I9").Value = 0.9630320723052701
$ws.Range("J9").Value = 0.9630320723052702
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.8212579999999999
$ws.Range("N9").Value = 2.463774
$ws.Range("O9").Value = 0.1112452839457698
$ws.Range("P9").Value = 0.1112452839457698
$ws.Range("Q9").Value = 1344.681319459246
$ws.Range("R9").Value = 12102.13187513322
$ws.Range("S9").Value = 0.1071327763324829
$ws.Range("T9").Value = 0.1071327763324829
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.50081933333334
$ws.Range("H10").Value = 52.502458
$ws.Range("I10").Value = 0.01029341242216722
$ws.Range("J10").Value = 0.01029341242216722
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.135193666666667
$ws.Range("N10").Value = 3.405581
$ws.Range("O10").Value = 0.153770120695047
$ws.Range("P10").Value = 0.153770120695047
$ws.Range("Q10").Value = 19.86681926867756
$ws.Range("R10").Value = 178.801373418098
$ws.Range("S10").Value = 0.001582819270520549
$ws.Range("T10").Value = 0.001582819270520549
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 17.50081933333334
$ws.Range("H11").Value = 52.502458
$ws.Range("I11").Value = 0.01029341242216722
$ws.Range("J11").Value = 0.01029341242216722
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.287366
$ws.Range("N11").Value = 6.862098
$ws.Range("O11").Value = 0.3098401235152652
$ws.Range("P11").Value = 0.3098401235152652
$ws.Range("Q11").Value = 40.03077911520934
$ws.Range("R11").Value = 360.277012036884
$ws.Range("S11").Value = 0.003189312176277856
$ws.Range("T11").Value = 0.003189312176277857
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 17.50081933333334
$ws.Range("H12").Value = 52.502458
$ws.Range("I12").Value = 0.01029341242216722
$ws.Range("J12").Value = 0.01029341242216722
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.138589666666667
$ws.Range("N12").Value = 9.415769000000001
$ws.Range("O12").Value = 0.425144471843918
$ws.Range("P12").Value = 0.425144471843918
$ws.Range("Q12").Value = 54.92789071780023
$ws.Range("R12").Value = 494.3510164602021
$ws.Range("S12").Value = 0.004376187387693906
$ws.Range("T12").Value = 0.004376187387693907
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 17.50081933333334
$ws.Range("H13").Value = 52.502458
$ws.Range("I13").Value = 0.01029341242216722
$ws.Range("J13").Value = 0.01029341242216722
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.8212579999999999
$ws.Range("N13").Value = 2.463774
$ws.Range("O13").Value = 0.1112452839457698
$ws.Range("P13").Value = 0.1112452839457698
$ws.Range("Q13").Value = 14.37268788405467
$ws.Range("R13").Value = 129.354190956492
$ws.Range("S13").Value = 0.001145093587674906
$ws.Range("T13").Value = 0.001145093587674907
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 25.39612333333333
$ws.Range("H14").Value = 76.18836999999999
$ws.Range("I14").Value = 0.01493717330687017
$ws.Range("J14").Value = 0.01493717330687017
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.135193666666667
$ws.Range("N14").Value = 3.405581
$ws.Range("O14").Value = 0.153770120695047
$ws.Range("P14").Value = 0.153770120695047
$ws.Range("Q14").Value = 28.82951836588555
$ws.Range("R14").Value = 259.46566529297
$ws.Range("S14").Value = 0.00229689094224026
$ws.Range("T14").Value = 0.00229689094224026
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 25.39612333333333
$ws.Range("H15").Value = 76.18836999999999
$ws.Range("I15").Value = 0.01493717330687017
$ws.Range("J15").Value = 0.01493717330687017
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.287366
$ws.Range("N15").Value = 6.862098
$ws.Range("O15").Value = 0.3098401235152652
$ws.Range("P15").Value = 0.3098401235152652
$ws.Range("Q15").Value = 58.09022904447333
$ws.Range("R15").Value = 522.8120614002599
$ws.Range("S15").Value = 0.004628135622369575
$ws.Range("T15").Value = 0.004628135622369575
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 25.39612333333333
$ws.Range("H16").Value = 76.18836999999999
$ws.Range("I16").Value = 0.01493717330687017
$ws.Range("J16").Value = 0.01493717330687017
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.138589666666667
$ws.Range("N16").Value = 9.415769000000001
$ws.Range("O16").Value = 0.425144471843918
$ws.Range("P16").Value = 0.425144471843918
$ws.Range("Q16").Value = 79.70801026739223
$ws.Range("R16").Value = 717.37209240653
$ws.Range("S16").Value = 0.006350456656390387
$ws.Range("T16").Value = 0.006350456656390388
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 25.39612333333333
$ws.Range("H17").Value = 76.18836999999999
$ws.Range("I17").Value = 0.01493717330687017
$ws.Range("J17").Value = 0.01493717330687017
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.8212579999999999
$ws.Range("N17").Value = 2.463774
$ws.Range("O17").Value = 0.1112452839457698
$ws.Range("P17").Value = 0.1112452839457698
$ws.Range("Q17").Value = 20.85676945648666
$ws.Range("R17").Value = 187.71092510838
$ws.Range("S17").Value = 0.001661690085869945
$ws.Range("T17").Value = 0.001661690085869945
